$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 7 (the old "totals" row becomes row 8)
$ws.Rows.Item(7).Insert()

# Fill in the new row 7 with the duplicated "proprietaire" contract data
$ws.Range("A7").Value = "040/SUP SUD"
$ws.Range("B7").Value = "Supervision"
$ws.Range("C7").Value = "BG1949"
$ws.Range("D7").Value = "Ahmed Test"
$ws.Range("E7").Value = "non"
$ws.Range("F7").Value = "mensuelle"
$ws.Range("G7").Value = 15
$ws.Range("H7").Value = "--"
$ws.Range("I7").Value = 15000
$ws.Range("J7").Value = "--"
$ws.Range("K7").Value = 2250
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 12750

# Recalculate the totals row (now row 8) to include the new contract's amounts
$ws.Range("I8").Value = 115000
$ws.Range("K8").Value = 3050
$ws.Range("M8").Value = 151950
